$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Rename template placeholder "startM" -> "abfahrtenM" everywhere
#    (3 occurrences, each inside "${startM}").
# ------------------------------------------------------------------
$d.Content.Find.Execute("startM", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "abfahrtenM", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark away from the "Alle anderen Termine"
#    heading paragraph to the end of the second "${abfahrtenM}"
#    paragraph (right after the closing "}", before the paragraph
#    mark).
# ------------------------------------------------------------------

# 2a. Remove the existing _GoBack bookmark (its old location).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2b. Locate the second "${abfahrtenM}" occurrence in document order.
$searchRange = $d.Content
$matchEnd = $null
$count = 0
while ($searchRange.Find.Execute("abfahrtenM") -eq $true) {
    $count = $count + 1
    if ($count -eq 2) {
        $matchEnd = $searchRange.End
        break
    }
    $searchRange.Collapse(0)
    $searchRange.End = $d.Content.End
}

# 2c. The run order is "${" "abfahrtenM" "}" - step over the closing
#     brace to land right after it, still inside the same paragraph.
$target = $matchEnd + 1

# 2d. Insert the new (collapsed) "_GoBack" bookmark at that position.
#     A bookmark added on a truly zero-width Range that sits exactly on
#     the paragraph mark is not placed reliably, so a temporary
#     placeholder character is inserted, bookmarked, and then removed
#     again - leaving the bookmark collapsed at the right spot.
$insertionPoint = $d.Range($target, $target)
$insertionPoint.InsertAfter("X")
$placeholder = $d.Range($target, $target + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$d.Range($target, $target + 1).Delete()
